$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix load-file SOI initiation: the points list for the last line (B8) was
# missing a point, so append it to the existing value.
$ws.Range("B8").Value = "Point_18 Point_20 Point_15"

# Reflect the resulting cursor/selection position (moved on to the next row).
$ws.Range("B10").Select()
